$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Addr="D2"; Value="66.534.97"}
    @{Addr="E2"; Value="  +0.42%  "}
    @{Addr="D3"; Value="3.217.55"}
    @{Addr="E3"; Value="  +1.15%  "}
    @{Addr="E4"; Value="  -0.02%  "}
    @{Addr="D5"; Value="609.33"}
    @{Addr="E5"; Value="  +2.39%  "}
    @{Addr="D6"; Value="158.12"}
    @{Addr="E6"; Value="  +2.81%  "}
    @{Addr="E7"; Value="  +0.00%  "}
    @{Addr="D8"; Value="3.216.75"}
    @{Addr="E8"; Value="  +1.17%  "}
    @{Addr="D9"; Value="0.551"}
    @{Addr="E9"; Value="  +0.41%  "}
    @{Addr="E10"; Value="  +0.77%  "}
    @{Addr="D11"; Value="5.69"}
    @{Addr="E11"; Value="  -4.24%  "}
    @{Addr="D12"; Value="0.503"}
    @{Addr="E12"; Value="  -2.55%  "}
    @{Addr="E13"; Value="  +0.73%  "}
    @{Addr="D14"; Value="38.74"}
    @{Addr="E14"; Value="  -0.76%  "}
    @{Addr="D15"; Value="3.746.86"}
    @{Addr="E15"; Value="  +1.12%  "}
    @{Addr="D16"; Value="66.638.98"}
    @{Addr="E16"; Value="  +0.63%  "}
    @{Addr="E17"; Value="  -0.99%  "}
    @{Addr="D18"; Value="3.218.32"}
    @{Addr="E18"; Value="  +0.89%  "}
    @{Addr="E19"; Value="  +1.31%  "}
    @{Addr="D20"; Value="506.16"}
    @{Addr="E20"; Value="  -1.55%  "}
    @{Addr="D21"; Value="15.16"}
    @{Addr="E21"; Value="  -1.27%  "}
    @{Addr="D22"; Value="0.733"}
    @{Addr="E22"; Value="  -0.57%  "}
    @{Addr="E23"; Value="  -0.50%  "}
    @{Addr="D24"; Value="14.61"}
    @{Addr="E24"; Value="  -2.17%  "}
    @{Addr="D25"; Value="84.88"}
    @{Addr="E25"; Value="  -0.80%  "}
    @{Addr="E26"; Value="  +0.43%  "}
    @{Addr="D27"; Value="3.00"}
    @{Addr="E27"; Value="  +0.10%  "}
    @{Addr="D28"; Value="9.12"}
    @{Addr="E28"; Value="  -1.14%  "}
    @{Addr="E29"; Value="  +1.19%  "}
    @{Addr="E30"; Value="  +36.47%  "}
    @{Addr="D31"; Value="2.94"}
    @{Addr="E31"; Value="  +0.93%  "}
    @{Addr="D32"; Value="7.02"}
    @{Addr="E32"; Value="  -1.32%  "}
    @{Addr="D33"; Value="28.12"}
    @{Addr="E33"; Value="  -0.42%  "}
    @{Addr="E34"; Value="  -0.06%  "}
    @{Addr="E35"; Value="  -4.04%  "}
    @{Addr="D36"; Value="6.48"}
    @{Addr="E36"; Value="  -0.36%  "}
    @{Addr="D37"; Value="502.44"}
    @{Addr="E37"; Value="  +0.13%  "}
    @{Addr="D38"; Value="55.41"}
    @{Addr="E38"; Value="  +1.06%  "}
    @{Addr="D39"; Value="0.0₃0770"}
    @{Addr="E39"; Value="  +14.11%  "}
    @{Addr="D40"; Value="0.131"}
    @{Addr="E40"; Value="  +5.53%  "}
    @{Addr="D41"; Value="0.0421"}
    @{Addr="E41"; Value="  -0.45%  "}
    @{Addr="D42"; Value="3.06"}
    @{Addr="E42"; Value="  +6.69%  "}
    @{Addr="D43"; Value="8.72"}
    @{Addr="E43"; Value="  -1.74%  "}
    @{Addr="D44"; Value="0.296"}
    @{Addr="E44"; Value="  -2.63%  "}
    @{Addr="E45"; Value="  +0.57%  "}
    @{Addr="D46"; Value="2.904.73"}
    @{Addr="E46"; Value="  -0.01%  "}
    @{Addr="E47"; Value="  -1.51%  "}
    @{Addr="E48"; Value="  +3.49%  "}
    @{Addr="E50"; Value="  -0.55%  "}
    @{Addr="D51"; Value="122.73"}
    @{Addr="E51"; Value="  -0.75%  "}
)

foreach ($u in $updates) {
    $ws.Range($u.Addr).Value = $u.Value
}
